{"js": "// 1) Delete the two obsolete table rows from the first table (HDS_new_pump):\n//      PUMP:HRD:0000  / [PUMP:TBV:1111]\n//      PUMP:HRD:3350  / [PUMP:TBV:1111]   (the duplicate-tag row right after the kept PUMP:HRD:3350 / [PUMP:HRS:3350] row)\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst firstTable = tables.items[0];\nfirstTable.rows.load(\"items\");\nawait context.sync();\n\nfor (const row of firstTable.rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nfor (const row of firstTable.rows.items) {\n  for (const cell of row.cells.items) {\n    cell.body.load(\"text\");\n  }\n}\nawait context.sync();\n\nconst rowsToDelete = [];\nfor (const row of firstTable.rows.items) {\n  const texts = row.cells.items.map((c) => c.body.text.trim());\n  const frontTag = texts[0];\n  const backTag = texts[1];\n  if (\n    frontTag === \"PUMP:HRD:0000\" ||\n    (frontTag === \"PUMP:HRD:3350\" && backTag === \"[PUMP:TBV:1111]\")\n  ) {\n    rowsToDelete.push(row);\n  }\n}\n\n// Delete from the bottom up so earlier indices stay valid, syncing between\n// each delete so the host re-resolves row positions.\nrowsToDelete.reverse();\nfor (const row of rowsToDelete) {\n  row.delete();\n  await context.sync();\n}\n\n// 2) Fix the stray space inserted after \"PUMP:\" (and the missing colon in\n//    \"PUMP HTR:200\") across the HTP_new_pump / HTR_new_pump tables.\nconst textFixes = [\n  [\"PUMP: HTP:200 \", \"PUMP:HTP:200 \"],\n  [\"PUMP: HTP:300 \", \"PUMP:HTP:300 \"],\n  [\"PUMP: HTP:400 \", \"PUMP:HTP:400 \"],\n  [\"PUMP: HTP:500 \", \"PUMP:HTP:500 \"],\n  [\"PUMP HTR:200 \", \"PUMP:HTR:200 \"],\n  [\"PUMP: HTR:300 \", \"PUMP:HTR:300 \"],\n  [\"PUMP: HTR:400 \", \"PUMP:HTR:400 \"],\n  [\"PUMP: HTR:500 \", \"PUMP:HTR:500 \"],\n  [\"PUMP: HTR:1100 \", \"PUMP:HTR:1100 \"],\n  [\"PUMP: HTR:1200 \", \"PUMP:HTR:1200 \"],\n  [\"PUMP: HTR:1300 \", \"PUMP:HTR:1300 \"],\n  [\"PUMP: HTR:1400 \", \"PUMP:HTR:1400 \"],\n  [\"PUMP: HTR:1500 \", \"PUMP:HTR:1500 \"],\n];\n\nfor (const [find, replace] of textFixes) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 3) Remove the space between the two bracketed tags in the\n//    SRS_BolusCalc_Pump_X04 table (BOLUS:SRS:2 row).\nconst prsTbd = context.document.body.search(\" [PUMP:PRS:1] [PUMP:TBD:1]\", {\n  matchCase: true,\n});\nprsTbd.load(\"items\");\nawait context.sync();\nfor (const item of prsTbd.items) {\n  item.insertText(\" [PUMP:PRS:1][PUMP:TBD:1]\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Delete the two obsolete table rows from the first table (HDS_new_pump):\n#      PUMP:HRD:0000  / [PUMP:TBV:1111]\n#      PUMP:HRD:3350  / [PUMP:TBV:1111]   (the duplicate-tag row right after the kept PUMP:HRD:3350 / [PUMP:HRS:3350] row)\n$t = $d.Tables.Item(1)\n\n$rowsToDelete = @()\nfor ($i = 1; $i -le $t.Rows.Count; $i++) {\n    # Cell().Range.Text carries a trailing cell-mark (CR + BEL); strip it before comparing.\n    $frontTag = $t.Cell($i, 1).Range.Text.TrimEnd([char]13, [char]7).Trim()\n    $backTag = $t.Cell($i, 2).Range.Text.TrimEnd([char]13, [char]7).Trim()\n    if (($frontTag -eq \"PUMP:HRD:0000\") -or (($frontTag -eq \"PUMP:HRD:3350\") -and ($backTag -eq \"[PUMP:TBV:1111]\"))) {\n        $rowsToDelete += $i\n    }\n}\n\n# Delete from the bottom up so earlier row indices stay valid.\n$rowsToDelete = $rowsToDelete | Sort-Object -Descending\nforeach ($rowIndex in $rowsToDelete) {\n    $t.Rows.Item($rowIndex).Delete()\n}\n\n# 2) Fix the stray space inserted after \"PUMP:\" (and the missing colon in\n#    \"PUMP HTR:200\") across the HTP_new_pump / HTR_new_pump tables.\n$textFixes = @(\n    @(\"PUMP: HTP:200 \", \"PUMP:HTP:200 \"),\n    @(\"PUMP: HTP:300 \", \"PUMP:HTP:300 \"),\n    @(\"PUMP: HTP:400 \", \"PUMP:HTP:400 \"),\n    @(\"PUMP: HTP:500 \", \"PUMP:HTP:500 \"),\n    @(\"PUMP HTR:200 \", \"PUMP:HTR:200 \"),\n    @(\"PUMP: HTR:300 \", \"PUMP:HTR:300 \"),\n    @(\"PUMP: HTR:400 \", \"PUMP:HTR:400 \"),\n    @(\"PUMP: HTR:500 \", \"PUMP:HTR:500 \"),\n    @(\"PUMP: HTR:1100 \", \"PUMP:HTR:1100 \"),\n    @(\"PUMP: HTR:1200 \", \"PUMP:HTR:1200 \"),\n    @(\"PUMP: HTR:1300 \", \"PUMP:HTR:1300 \"),\n    @(\"PUMP: HTR:1400 \", \"PUMP:HTR:1400 \"),\n    @(\"PUMP: HTR:1500 \", \"PUMP:HTR:1500 \")\n)\n\nforeach ($pair in $textFixes) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\n# 3) Remove the space between the two bracketed tags in the\n#    SRS_BolusCalc_Pump_X04 table (BOLUS:SRS:2 row).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" [PUMP:PRS:1] [PUMP:TBD:1]\"\n$find.Replacement.Text = \" [PUMP:PRS:1][PUMP:TBD:1]\"\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
